$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 to make room for the "Not applicable" / -1 entry.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row.
$ws.Range("A2").Value = -1
$ws.Range("B2").Value = "Not applicable"

# Update the active selection to match the target state.
$ws.Range("A3").Select()
